$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 data values ---
# New strings are introduced in alphabetical order so the shared-string
# table they get appended to ends up sorted the same way the source
# workbook has it.
$ws.Range("D2").Value = "Jharsuguda"
$ws.Range("E2").Value = "Jharsuguda"
$ws.Range("F2").Value = "Laikera"
$ws.Range("G2").Value = "LAIKERA SCS"
$ws.Range("A2").Value = "sasi"
$ws.Range("B2").Value = "sasi11@gmail.com"
$ws.Range("C2").Value = "Siba@123"

$ws.Range("I2").Value = "Rama"
$ws.Range("O2").Value = "CUTTACK"
$ws.Range("P2").Value = "Jagatsinghpur"
$ws.Range("Q2").Value = "Niali"
$ws.Range("T2").Value = "payal@oasys.com"

# R2 is removed entirely (no submit status recorded any more)
$ws.Range("R2").ClearContents() | Out-Null

# C2 becomes a hyperlink (mirrors B2's hyperlinked e-mail pattern)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:sasi11@gmail.com") | Out-Null

# --- Column width tweaks (column F widened, column G gets an explicit width) ---
$ws.Range("F1").ColumnWidth = 18
$ws.Range("G1").ColumnWidth = 17.666666666666668

# --- Selection / view state ---
$ws.Range("R5").Select() | Out-Null
